$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2 (the single existing data row) into row 3, keeping the
# exact same text-stored values/formatting as row 2 (copy + paste-values
# avoids Excel re-interpreting the numeric-looking strings as real numbers).
$ws.Range("A2:K2").Copy()
$ws.Range("A3").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# The new row's batsman name uses a plain trailing space rather than the
# non-breaking space present in row 2.
$ws.Cells.Item(3, 6).Value = "Dhawal Kulkarni "
